$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update Change Type for Part Number / "Show the Part Details based on Jobcard" row from Small to Medium
$ws.Range("C8").Value = "Medium"

# Fill in Man Hours column (G) with estimated hours for each change row
$ws.Range("G2").Value = 2
$ws.Range("G3").Value = 2
$ws.Range("G4").Value = 0.5
$ws.Range("G5").Value = 3
$ws.Range("G6").Value = 6
$ws.Range("G7").Value = 1
$ws.Range("G8").Value = 8
$ws.Range("G9").Value = 1
$ws.Range("G10").Value = 1
$ws.Range("G11").Value = 1
$ws.Range("G12").Value = 8
$ws.Range("G13").Value = 2
$ws.Range("G14").Value = 1
$ws.Range("G15").Value = 6
$ws.Range("G16").Value = 8
$ws.Range("G17").Value = 6
$ws.Range("G18").Value = 4
$ws.Range("G19").Value = 16
$ws.Range("G20").Value = 24
$ws.Range("G21").Value = 12
$ws.Range("G22").Value = 4
$ws.Range("G23").Value = 16

# Update the view: scroll so row 7 is at top-left, and select the full Man Hours data range
$ws.Range("G2:G23").Select()
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Application.ActiveWindow.ScrollColumn = 1
